$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Chapter 2 ("2-Marco teorico", row 3) is now read for Martin (column C):
# previously marked "En proceso", now marked "Leido".
$ws.Range("C3").Value = "Leido"

# The "En proceso" mark moves down to row 5 ("4-Metodologia de Desarrollo"),
# which previously had no value in column C.
$ws.Range("C5").Value = "En proceso"

# Update the active selection to reflect the last edited cell.
$ws.Range("C5").Select()
